$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 196, pushing existing rows 196:232 down to 197:233.
$ws.Rows.Item(196).EntireRow.Insert()

# Populate the newly inserted row 196 with the new price-observation record.
$ws.Range("A196").Value = 11
$ws.Range("B196").Value = "Vega Monumental Concepción"
$ws.Range("C196").Value = "Bíobío"
$ws.Range("D196").Value = 44551
$ws.Range("E196").Value = 8
$ws.Range("F196").Value = "Fruta"
$ws.Range("G196").Value = 100101
$ws.Range("H196").Value = "Berries"
$ws.Range("I196").Value = 100112025
$ws.Range("J196").Value = "Frutilla"
$ws.Range("K196").Value = "Sin especificar"
$ws.Range("L196").Value = "Primera"
$ws.Range("M196").Value = 270
$ws.Range("N196").Value = 7000
$ws.Range("O196").Value = 7500
$ws.Range("P196").Value = 7222
$ws.Range("Q196").Value = "$/bandeja 7 kilos"
$ws.Range("R196").Value = "Región Metropolitana"
$ws.Range("S196").Value = 1032
$ws.Range("T196").Value = 7
